# Update the "Förändrad" (Changed) date column (C) from 2023-09-06 (45175)
# to 2023-09-08 (45177) for every data row in the sheet (header is row 1,
# data starts on row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row by walking up from the bottom of column A
# (xlUp = -4162), so the script keeps working if the sheet grows/shrinks.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45177
}
